# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Writes the newly-computed "K" values (column G) back into the sheet for
# rows 2-33 of the existing save-data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 6
    4  = 4
    5  = 6
    6  = 4
    7  = 4
    8  = 7
    9  = 5
    10 = 4
    11 = 6
    12 = 5
    13 = 5
    14 = 7
    15 = 6
    16 = 7
    17 = 7
    18 = 6
    19 = 5
    20 = 9
    21 = 7
    22 = 2
    23 = 6
    24 = 4
    25 = 4
    26 = 9
    27 = 9
    28 = 5
    29 = 6
    30 = 7
    31 = 1
    32 = 6
    33 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
